# @dev add axios work
# Append a new worklist row (row 10) describing the axios "add shopping
# config" task, matching the existing table layout (#, description, file, status).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = "axios 增加购物配置"
$ws.Cells.Item(10, 3).Value = "AddConfig"
$ws.Cells.Item(10, 4).Value = "未做"

$ws.Range("D10").Select()
